# Update the "攻略人数" (F column) figures across all four sheets to match
# the newly published gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2551
$ws.Range("F3").Value = 550
$ws.Range("F4").Value = 453
$ws.Range("F5").Value = 284
$ws.Range("F7").Value = 449
$ws.Range("F8").Value = 1184
$ws.Range("F9").Value = 534
$ws.Range("F10").Value = 292
$ws.Range("F11").Value = 113
$ws.Range("F12").Value = 342
$ws.Range("F13").Value = 5497
$ws.Range("F14").Value = 48
$ws.Range("F15").Value = 1661
$ws.Range("F16").Value = 3985
$ws.Range("F17").Value = 404
$ws.Range("F20").Value = 4542
$ws.Range("F21").Value = 5947
$ws.Range("F23").Value = 1016
$ws.Range("F24").Value = 654
$ws.Range("F25").Value = 3666
$ws.Range("F26").Value = 469
$ws.Range("F27").Value = 64
$ws.Range("F28").Value = 180
$ws.Range("F29").Value = 117
$ws.Range("F30").Value = 961
$ws.Range("F31").Value = 1356
$ws.Range("F32").Value = 445
$ws.Range("F33").Value = 508
$ws.Range("F35").Value = 188
$ws.Range("F36").Value = 1640
$ws.Range("F37").Value = 156
$ws.Range("F38").Value = 7
$ws.Range("F39").Value = 1073
$ws.Range("F41").Value = 1338
$ws.Range("F42").Value = 601
$ws.Range("F45").Value = 3259
$ws.Range("F47").Value = 258
$ws.Range("F49").Value = 3854

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 5
$ws.Range("F5").Value = 1170
$ws.Range("F22").Value = 65

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 3670

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 2551
$ws.Range("F4").Value = 550
$ws.Range("F5").Value = 453
$ws.Range("F6").Value = 284
$ws.Range("F7").Value = 1170
$ws.Range("F9").Value = 449
$ws.Range("F10").Value = 1184
$ws.Range("F11").Value = 534
$ws.Range("F12").Value = 292
$ws.Range("F13").Value = 113
$ws.Range("F14").Value = 342
$ws.Range("F15").Value = 5497
$ws.Range("F17").Value = 1661
$ws.Range("F18").Value = 4542
$ws.Range("F19").Value = 5947
$ws.Range("F21").Value = 1016
$ws.Range("F22").Value = 654
$ws.Range("F23").Value = 3666
$ws.Range("F24").Value = 469
$ws.Range("F25").Value = 64
$ws.Range("F26").Value = 180
$ws.Range("F27").Value = 117
$ws.Range("F28").Value = 1356
$ws.Range("F29").Value = 445
$ws.Range("F30").Value = 508
$ws.Range("F33").Value = 188
$ws.Range("F34").Value = 1640
$ws.Range("F36").Value = 1073
$ws.Range("F38").Value = 601
$ws.Range("F42").Value = 65
$ws.Range("F43").Value = 3259
$ws.Range("F46").Value = 258
$ws.Range("F49").Value = 3854
